# Insert a new weekly price-record row for "Apio" (Feria Lagunitas de Puerto
# Montt) at sheet row 72, pushing the existing rows 72..174 down to 73..175.
#
# The new row carries the same Mercado/Región/Codreg/Categoría/Variedad/
# Calidad/Unidad/Origen/Kg-o-Unidades/Clasificación values as the record that
# used to sit at row 72 (it's the same market "Primera" series), but with a
# newer Fecha and updated Volumen/Precio columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 72:174 down to 73:175, leaving a blank row 72 behind (Excel
# automatically carries over row/column formatting, e.g. the date style on
# column D).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record.
$ws.Range("A72").Value2 = 4
$ws.Range("B72").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value2 = "Los Lagos"
$ws.Range("D72").Value2 = 44546
$ws.Range("E72").Value2 = 10
$ws.Range("F72").Value2 = 100112017
$ws.Range("G72").Value2 = "Apio"
$ws.Range("H72").Value2 = "Americana (o)"
$ws.Range("I72").Value2 = "Primera"
$ws.Range("J72").Value2 = 25
$ws.Range("K72").Value2 = 12000
$ws.Range("L72").Value2 = 12000
$ws.Range("M72").Value2 = 12000
$ws.Range("N72").Value2 = "`$/docena de matas"
$ws.Range("O72").Value2 = "Región de Coquimbo"
$ws.Range("P72").Value2 = 2000
$ws.Range("Q72").Value2 = 6
$ws.Range("R72").Value2 = "Hortaliza"
